$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "267.80"
$c.Style = "Normal"

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "22.90"
$c.Style = "Normal"

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "0.06210"
$c.Style = "Normal"

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "3.590"
$c.Style = "Normal"

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "6.681"
$c.Style = "Normal"

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.8326"
$c.Style = "Normal"

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.01363"
$c.Style = "Normal"

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.1599"
$c.Style = "Normal"

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.08251"
$c.Style = "Normal"

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.03398"
$c.Style = "Normal"

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.09278"
$c.Style = "Normal"

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "3.919"
$c.Style = "Normal"

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.001716"
$c.Style = "Normal"

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.04844"
$c.Style = "Normal"

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.006262"
$c.Style = "Normal"

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "0.005382"
$c.Style = "Normal"

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "3.756"
$c.Style = "Normal"

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "2.328"
$c.Style = "Normal"

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "0.1221"
$c.Style = "Normal"

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0.0002682"
$c.Style = "Normal"

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.04650"
$c.Style = "Normal"

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.006872"
$c.Style = "Normal"

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.1156"
$c.Style = "Normal"

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.003349"
$c.Style = "Normal"

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.01212"
$c.Style = "Normal"

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.00006220"
$c.Style = "Normal"

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.7000"
$c.Style = "Normal"

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.1672"
$c.Style = "Normal"

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.00002100"
$c.Style = "Normal"

